$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text
# (e.g. "41.00", "2.377.45"). Force those cells to the Text number
# format before assigning so Excel does not auto-convert them to
# numbers and silently drop formatting such as trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.315.41'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.369.96'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.55'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.65'
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.636'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.00'
$ws.Range('E10').Value = '  -4.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0918'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.49'
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.109'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.983'
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.732.11'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.41'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.368.50'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.307.07'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.26'
$ws.Range('E19').Value = '  +14.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.28'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000106'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.66'
$ws.Range('E22').Value = '  +3.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.37'
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.21'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.49'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.17'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.30'
$ws.Range('E29').Value = '  -1.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.44'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0947'
$ws.Range('E31').Value = '  -2.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.10'
$ws.Range('E32').Value = '  -4.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '169.08'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.81'
$ws.Range('E34').Value = '  -5.68%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('E36').Value = '  -2.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.71'
$ws.Range('E37').Value = '  -5.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.91'
$ws.Range('E38').Value = '  +10.27%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.97'
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0355'
$ws.Range('E41').Value = '  -3.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.59'
$ws.Range('E42').Value = '  -5.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.20'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.875.98'
$ws.Range('E44').Value = '  +14.31%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.227'
$ws.Range('E45').Value = '  -5.14%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.88'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.95'
$ws.Range('E48').Value = '  +3.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '83.44'
$ws.Range('E49').Value = '  +5.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '111.66'
$ws.Range('E50').Value = '  -3.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.21'
$ws.Range('E51').Value = '  -1.77%  '
